$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Withdrawn" counts (value = 1) that were left blank for
# Army, Navy, and Marine Corps sections.
$ws.Range("B26").Value = 1
$ws.Range("B32").Value = 1
$ws.Range("B37").Value = 1

# Reset the active selection back to the top-left cell (matches the
# saved file no longer pinning the cursor at A16).
$ws.Range("A1").Select()
